$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.992.13'
$ws.Range("E2").Value = '  -4.24%  '
$ws.Range("D3").Value = '2.229.92'
$ws.Range("E3").Value = '  -6.38%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.60'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -5.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '80.13'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -8.77%  '
$ws.Range("E7").Value = '  -4.19%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.459'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -6.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0770'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '27.85'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -10.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.83'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -13.76%  '
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").Value = '2.566.41'
$ws.Range("E14").Value = '  -6.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.09'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -7.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.08'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.95%  '
$ws.Range("D17").Value = '2.235.85'
$ws.Range("E17").Value = '  -6.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.715'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -5.97%  '
$ws.Range("D19").Value = '38.916.20'
$ws.Range("E19").Value = '  -4.23%  '
$ws.Range("D20").Value = '0.0₃0856'
$ws.Range("E20").Value = '  -6.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.72'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -7.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.80'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.79'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -9.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '224.74'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.93%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -10.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.71'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.10%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.17'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.89'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.73'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.09'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.76'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -9.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.34'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0683'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.58%  '
$ws.Range("E37").Value = '  -4.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.67'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0948'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.66%  '
$ws.Range("E40").Value = '  -7.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.46'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -9.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.64'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D43").Value = '1.905.16'
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("E44").Value = '  -8.89%  '
$ws.Range("E45").Value = '  -6.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.18'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -8.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.99'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.51'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -10.58%  '
$ws.Range("D49").Value = '2.435.64'
$ws.Range("E49").Value = '  -6.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '87.39'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.84'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -7.24%  '
